{"js": "// Replace each three-digit division equation/answer string in the\n// document's table cells with its updated value, per the commit diff.\n// Each \"before\" string is unique in the document, so a direct\n// search-and-replace (first match) is safe for every entry.\nconst replacements = [\n  [\"728\u00f75=145, 3\", \"471\u00f75=94, 1\"],\n  [\"867\u00f78=108, 3\", \"632\u00f74=158, 0\"],\n  [\"243\u00f74=60, 3\", \"443\u00f79=49, 2\"],\n  [\"331\u00f72=165, 1\", \"571\u00f78=71, 3\"],\n  [\"731\u00f78=91, 3\", \"246\u00f76=41, 0\"],\n  [\"737\u00f76=122, 5\", \"722\u00f72=361, 0\"],\n  [\"186\u00f76=31, 0\", \"549\u00f75=109, 4\"],\n  [\"866\u00f74=216, 2\", \"226\u00f79=25, 1\"],\n  [\"648\u00f76=108, 0\", \"189\u00f78=23, 5\"],\n  [\"542\u00f77=77, 3\", \"251\u00f79=27, 8\"],\n  [\"381\u00f79=42, 3\", \"970\u00f78=121, 2\"],\n  [\"551\u00f79=61, 2\", \"441\u00f77=63, 0\"],\n  [\"641\u00f77=91, 4\", \"363\u00f78=45, 3\"],\n  [\"939\u00f78=117, 3\", \"646\u00f72=323, 0\"],\n  [\"319\u00f79=35, 4\", \"950\u00f76=158, 2\"],\n  [\"927\u00f79=103, 0\", \"514\u00f78=64, 2\"],\n  [\"834\u00f78=104, 2\", \"469\u00f74=117, 1\"],\n  [\"749\u00f72=374, 1\", \"846\u00f77=120, 6\"],\n  [\"570\u00f79=63, 3\", \"806\u00f73=268, 2\"],\n  [\"845\u00f77=120, 5\", \"668\u00f72=334, 0\"],\n  [\"153\u00f77=21, 6\", \"187\u00f74=46, 3\"],\n  [\"651\u00f78=81, 3\", \"455\u00f73=151, 2\"],\n  [\"592\u00f72=296, 0\", \"298\u00f76=49, 4\"],\n  [\"124\u00f77=17, 5\", \"930\u00f75=186, 0\"],\n  [\"709\u00f78=88, 5\", \"668\u00f78=83, 4\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  // Replace only the first occurrence (old strings are unique in doc).\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace each three-digit division equation/answer string in the\n# document's table cells with its updated value, per the commit diff.\n# Each \"before\" string is unique in the document, so Find/Replace\n# (wdReplaceOne) targets exactly the one cell that should change.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"728\u00f75=145, 3\"; New = \"471\u00f75=94, 1\" }\n    @{ Old = \"867\u00f78=108, 3\"; New = \"632\u00f74=158, 0\" }\n    @{ Old = \"243\u00f74=60, 3\"; New = \"443\u00f79=49, 2\" }\n    @{ Old = \"331\u00f72=165, 1\"; New = \"571\u00f78=71, 3\" }\n    @{ Old = \"731\u00f78=91, 3\"; New = \"246\u00f76=41, 0\" }\n    @{ Old = \"737\u00f76=122, 5\"; New = \"722\u00f72=361, 0\" }\n    @{ Old = \"186\u00f76=31, 0\"; New = \"549\u00f75=109, 4\" }\n    @{ Old = \"866\u00f74=216, 2\"; New = \"226\u00f79=25, 1\" }\n    @{ Old = \"648\u00f76=108, 0\"; New = \"189\u00f78=23, 5\" }\n    @{ Old = \"542\u00f77=77, 3\"; New = \"251\u00f79=27, 8\" }\n    @{ Old = \"381\u00f79=42, 3\"; New = \"970\u00f78=121, 2\" }\n    @{ Old = \"551\u00f79=61, 2\"; New = \"441\u00f77=63, 0\" }\n    @{ Old = \"641\u00f77=91, 4\"; New = \"363\u00f78=45, 3\" }\n    @{ Old = \"939\u00f78=117, 3\"; New = \"646\u00f72=323, 0\" }\n    @{ Old = \"319\u00f79=35, 4\"; New = \"950\u00f76=158, 2\" }\n    @{ Old = \"927\u00f79=103, 0\"; New = \"514\u00f78=64, 2\" }\n    @{ Old = \"834\u00f78=104, 2\"; New = \"469\u00f74=117, 1\" }\n    @{ Old = \"749\u00f72=374, 1\"; New = \"846\u00f77=120, 6\" }\n    @{ Old = \"570\u00f79=63, 3\"; New = \"806\u00f73=268, 2\" }\n    @{ Old = \"845\u00f77=120, 5\"; New = \"668\u00f72=334, 0\" }\n    @{ Old = \"153\u00f77=21, 6\"; New = \"187\u00f74=46, 3\" }\n    @{ Old = \"651\u00f78=81, 3\"; New = \"455\u00f73=151, 2\" }\n    @{ Old = \"592\u00f72=296, 0\"; New = \"298\u00f76=49, 4\" }\n    @{ Old = \"124\u00f77=17, 5\"; New = \"930\u00f75=186, 0\" }\n    @{ Old = \"709\u00f78=88, 5\"; New = \"668\u00f78=83, 4\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $found = $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 1)\n    if (-not $found) {\n        throw \"Text not found: $($pair.Old)\"\n    }\n}\n"}
